# Insert two new weekly rows ("Primera" and "Segunda" quality records for
# 2023-11-XX / serial date 45265) at the top of the Brócoli data block
# (row 1206), pushing the existing 1206:1297 data down to 1208:1299.
# This grows the used range from A1:R1297 to A1:R1299.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 1206 (shifts old rows 1206..1297 -> 1208..1299).
$ws.Rows.Item(1206).Insert()
$ws.Rows.Item(1206).Insert()

# New row 1206: "Primera" quality record.
$ws.Cells.Item(1206, 1).Value = 3
$ws.Cells.Item(1206, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1206, 3).Value = "Coquimbo"
$ws.Cells.Item(1206, 4).Value = 45265
$ws.Cells.Item(1206, 5).Value = 5
$ws.Cells.Item(1206, 6).Value = 100112023
$ws.Cells.Item(1206, 7).Value = "Brócoli"
$ws.Cells.Item(1206, 8).Value = "Sin especificar"
$ws.Cells.Item(1206, 9).Value = "Primera"
$ws.Cells.Item(1206, 10).Value = 2600
$ws.Cells.Item(1206, 11).Value = 800
$ws.Cells.Item(1206, 12).Value = 900
$ws.Cells.Item(1206, 13).Value = 862
$ws.Cells.Item(1206, 14).Value = "$/unidad"
$ws.Cells.Item(1206, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1206, 16).Value = 862
$ws.Cells.Item(1206, 17).Value = 1
$ws.Cells.Item(1206, 18).Value = "Hortaliza"

# New row 1207: "Segunda" quality record.
$ws.Cells.Item(1207, 1).Value = 3
$ws.Cells.Item(1207, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1207, 3).Value = "Coquimbo"
$ws.Cells.Item(1207, 4).Value = 45265
$ws.Cells.Item(1207, 5).Value = 5
$ws.Cells.Item(1207, 6).Value = 100112023
$ws.Cells.Item(1207, 7).Value = "Brócoli"
$ws.Cells.Item(1207, 8).Value = "Sin especificar"
$ws.Cells.Item(1207, 9).Value = "Segunda"
$ws.Cells.Item(1207, 10).Value = 1200
$ws.Cells.Item(1207, 11).Value = 700
$ws.Cells.Item(1207, 12).Value = 700
$ws.Cells.Item(1207, 13).Value = 700
$ws.Cells.Item(1207, 14).Value = "$/unidad"
$ws.Cells.Item(1207, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1207, 16).Value = 700
$ws.Cells.Item(1207, 17).Value = 1
$ws.Cells.Item(1207, 18).Value = "Hortaliza"
